# Add the new sheet "工作表1" as the last sheet in the workbook and
# populate it with the employee name / learned-skills / strength table.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "工作表1"

$data = @(
    @("姓名", "學過", "擅長"),
    @("陳泓毓", "Py,Ts", "Py,Ts"),
    @("林志明", "Js,Java", "Java"),
    @("張雅婷", "SQL,C++", "SQL"),
    @("李宗翰", "Py,Java", "Py"),
    @("王淑芬", "Ts,SQL", "Ts"),
    @("黃俊傑", "C++,Py", "Py"),
    @("吳佩珊", "Java,Ts", "Java"),
    @("劉家豪", "SQL,Js", "SQL"),
    @("蔡依林", "Py,Ts,Java", "Py,Java"),
    @("許文彬", "Ts,Js", "Ts"),
    @("鄭美玲", "Java,SQL", "SQL"),
    @("謝志豪", "C++,Py,Ts", "Py"),
    @("洪嘉欣", "Js,Py", "Js"),
    @("郭俊宏", "SQL,Java", "SQL"),
    @("邱雅雯", "Ts,C++", "Ts"),
    @("陳建宏", "Py,Js", "Py"),
    @("林佳慧", "Java,Py", "Java"),
    @("張志強", "SQL,Ts", "SQL"),
    @("李佩玲", "C++,Js", "Js"),
    @("王建銘", "Py,SQL", "Py"),
    @("黃雅慧", "Ts,Java", "Java"),
    @("吳宗憲", "Js,C++", "C++"),
    @("劉俊宏", "SQL,Py", "SQL"),
    @("蔡佳玲", "Java,Ts,Py", "Java,Py"),
    @("許志明", "C++,SQL", "SQL"),
    @("鄭雅婷", "Py,Js,Ts", "Py"),
    @("謝宗翰", "Ts,SQL,Java", "SQL"),
    @("洪佩珊", "Js,Py,SQL", "Js"),
    @("郭怡君", "Java,C++", "Java"),
    @("邱俊傑", "SQL,Ts,Py", "Ts")
)

# Column A (names) first ...
for ($i = 1; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $data[$i][0]
}
# ... then the header row ...
$ws.Cells.Item(1, 1).Value = $data[0][0]
$ws.Cells.Item(1, 2).Value = $data[0][1]
$ws.Cells.Item(1, 3).Value = $data[0][2]
# ... then column B ...
for ($i = 1; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $data[$i][1]
}
# ... then column C.
for ($i = 1; $i -lt $data.Count; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $data[$i][2]
}

$ws.Columns.Item(2).ColumnWidth = 14.1
$ws.Columns.Item(3).ColumnWidth = 19.25

$ws.Select() | Out-Null
$ws.Range("C17").Select() | Out-Null
